$d = $word.ActiveDocument

# --- Part 1: merge the "is the degrees latitude." runs into a single run,
# removing the gramStart/gramEnd proofErr markers that wrapped "latitude."
$d.Content.Find.Execute("is the degrees latitude. ", $false, $false, $false, $false, $false, $true, 1, $false, "is the degrees latitude. ", 2) | Out-Null

# --- Part 2: append two new paragraphs at the end of the document:
#   1) an empty paragraph
#   2) a paragraph with the new sentence about sensitivity tests
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newLast = $d.Paragraphs.Last
$newRange = $newLast.Range
$emptyParaPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($emptyParaPackage)

$finalPara = $d.Paragraphs.Last
$finalRange = $finalPara.Range
$newText = "We also conduct a suite of sensitivity tests to provide additional constraints on the error of the optimized emissions, which are summarized in section 2.6."
$insertStart = $finalRange.Start
$finalRange.InsertBefore($newText)

$typedRange = $d.Range($insertStart, $insertStart + $newText.Length)
$typedRange.Font.Name = "Times"
